# Apply the two kinds of edits described by the diff:
#  1. Column C ("Förändrad") date serial changes from 45184 to 45186 for every data row.
#  2. Every HYPERLINK(...) formula in columns S:Y gets a second argument equal to the
#     "Beteckning" value found in column A of the same row (the friendly link text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the full extent of the data.
$lastRow = $ws.UsedRange.Rows.Count - 1   # UsedRange is 0-based here (row 0 exists and is blank)
$firstDataRow = 2
$lastDataRow = $lastRow

# --- 1. Bump the "Förändrad" date (column C) for every data row at once -----------------
$ws.Range("C$firstDataRow`:C$lastDataRow").Value = 45186

# --- 2. Add the display-text argument to every HYPERLINK formula in S:Y -----------------
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {

    $label = $ws.Range("A$row").Value2
    if (-not $label) { continue }

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range("$col$row")

        if (-not $cell.HasFormula) { continue }

        $formula = $cell.Formula

        # Only touch plain single-argument HYPERLINK("...") formulas; skip ones that
        # already carry a friendly-name argument.
        if ($formula -match '^=HYPERLINK\(\s*"[^"]*"\s*\)$') {
            $newFormula = $formula -replace '\)\s*$', (', "' + $label + '")')
            $cell.Formula = $newFormula
        }
    }
}
